# Apply Feb 3 2024 cryptos list update (prices / 1h volume changes, and
# a row-41/42 coin swap) as captured by the upstream GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.028.12"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Formula = "'300.06"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("D6").Formula = "'97.88"
$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("E7").Value = "  +3.39%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").Formula = "'36.26"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("D11").Formula = "'0.0792"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").Formula = "'17.81"
$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").Formula = "'6.90"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").Value = "2.660.01"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "2.280.38"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").Value = "42.910.78"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Formula = "'12.99"
$ws.Range("E19").Value = "  +3.32%  "

$ws.Range("D20").Value = "0.0₃0912"
$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").Formula = "'68.25"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("D23").Formula = "'237.77"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("E24").Value = "  -1.19%  "

$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").Formula = "'2.43"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").Formula = "'25.00"

$ws.Range("D29").Formula = "'2.06"
$ws.Range("E29").Value = "  -12.53%  "

$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("D31").Formula = "'163.43"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("D33").Formula = "'0.999"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Formula = "'5.12"
$ws.Range("E34").Value = "  +2.11%  "

$ws.Range("D35").Formula = "'18.19"
$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("E36").Value = "  +3.44%  "

$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").Formula = "'0.0697"
$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("E39").Value = "  +0.81%  "

$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Formula = "'2.78"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Formula = "'0.111"
$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("D43").Value = "2.013.73"
$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("D44").Formula = "'0.0287"
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("E46").Value = "  +1.27%  "

$ws.Range("D47").Formula = "'17.52"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").Formula = "'2.84"
$ws.Range("E48").Value = "  -1.98%  "

$ws.Range("D49").Formula = "'54.30"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").Value = "2.530.45"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").Formula = "'1.53"
$ws.Range("E51").Value = "  -0.67%  "
